$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Capture the existing "Description" / "Purpose" / "Copyright" / "Immutable"
# rows (currently rows 11-14) before they get shifted down to make room for
# the new "Contact" (duplicate) and "Jurisdiction" rows.
$a11 = $ws.Cells.Item(11, 1).Value2
$b11 = $ws.Cells.Item(11, 2).Value2
$a12 = $ws.Cells.Item(12, 1).Value2
$b12 = $ws.Cells.Item(12, 2).Value2
$a13 = $ws.Cells.Item(13, 1).Value2
$b13 = $ws.Cells.Item(13, 2).Value2
$a14 = $ws.Cells.Item(14, 1).Value2
$b14 = $ws.Cells.Item(14, 2).Value2

# Extend the shared row formatting (border/alignment style used by every
# data row) down through the two new rows so nothing ends up unstyled.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B16").PasteSpecial(-4122)

# Updated Date value
$ws.Range("B8").Value = "2021-12-22T21:26:07+01:00"

# Updated Publisher value
$ws.Range("B9").Value = "Forschungsgruppe Digital Health"

# New row 11: repeated Contact / "No display for ContactDetail" entry
$ws.Cells.Item(11, 1).Value = "Contact"
$ws.Cells.Item(11, 2).Value = "No display for ContactDetail"

# New row 12: Jurisdiction / Germany
$ws.Cells.Item(12, 1).Value = "Jurisdiction"
$ws.Cells.Item(12, 2).Value = "Germany"

# Shifted-down rows 13-16 (previously rows 11-14)
$ws.Cells.Item(13, 1).Value = $a11
$ws.Cells.Item(13, 2).Value = $b11
$ws.Cells.Item(14, 1).Value = $a12
$ws.Cells.Item(14, 2).Value = $b12
$ws.Cells.Item(15, 1).Value = $a13
$ws.Cells.Item(15, 2).Value = $b13
$ws.Cells.Item(16, 1).Value = $a14
$ws.Cells.Item(16, 2).Value = $b14
